$d = $word.ActiveDocument

# --- Change 1: prepend the markdown code-fence marker + a space before
#     the "import platform" run in the FirstParagraph paragraph, as two
#     distinct runs: "’’’{python}" and " ".
$marker = "’’’{python}"
$p = $d.Paragraphs.Item(4)
$r = $p.Range
$ins = $d.Range($r.Start, $r.Start)
$ins.InsertBefore($marker + " ")

# Split the newly-inserted text "’’’{python} " away from the following
# "import platform" run (it's currently one merged run) by forcing a
# run boundary at the marker/space seam without introducing any
# lingering run-formatting (self-assigning FormattedText splits runs
# cleanly with no leftover <w:rPr/>).
$markerRange = $d.Range($r.Start, $r.Start + $marker.Length)
$markerRange.FormattedText = $markerRange.FormattedText

# --- Change 2: drop the trailing ``` fence from the end of the last
#     paragraph's final run.
$backtick = [string][char]96
$fence = $backtick + $backtick + $backtick
$old = "2), 2))" + $fence
$new = "2), 2))"
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
